# Finalizing BOM: update the two input cells that drive the torque calc
# (Approximate Mass C21 and the length/radius input C25). The dependent
# formulas in C23 (Acceleration Force), C26 and C27 (Torque) recalculate
# automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 50
$ws.Range("C25").Value = 0.13
